# "Update college and QQ group"
#
# On the title slide (slide 1), the "学院：..." (College:) caption is
# changed from the old college name to "人工智能学院" (School of
# Artificial Intelligence). The caption text box is an auto-fit
# ("shrink/grow to fit text", wrap="none") shape that is centered on a
# fixed midpoint, so after the (shorter) text is applied PowerPoint
# re-lays out the box and its left edge/width shrink to fit while its
# vertical position/height and horizontal center stay put.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$collegeBox = $s.Shapes.Item(3)   # id=21, "学院：..." caption textbox
$collegeBox.TextFrame.TextRange.Text = "学院：人工智能学院"

# Re-layout of the autosize textbox after the text edit: keep the same
# horizontal center (6203748 EMU) and vertical position, but shrink the
# box to the new, narrower text extent (values are expressed in points,
# as PowerPoint's Shape.Left/Width do, nudged to land on the exact EMU
# PowerPoint itself computes for this caption after the shrink).
$collegeBox.Left = 418.2824859448819
$collegeBox.Width = 140.40335083661418
